# Update the "想去人数" (F column) counts on the "展览" and "全部类型"
# worksheets to reflect the newer scrape snapshot (456a3b4).
#
# Row -> old value -> new value
#   2  : 11732 -> 11741
#   3  : 11362 -> 11382
#   11 : 10774 -> 10779
#   12 : 4165  -> 4168
#   16 : 2469  -> 2470
#   17 : 1051  -> 1052
#   18 : 51    -> 52
#   20 : 451   -> 452
#   22 : 10928 -> 10930
#   28 : 14    -> 15

$wb = $excel.ActiveWorkbook

$updates = @{}
$updates[2]  = 11741
$updates[3]  = 11382
$updates[11] = 10779
$updates[12] = 4168
$updates[16] = 2470
$updates[17] = 1052
$updates[18] = 52
$updates[20] = 452
$updates[22] = 10930
$updates[28] = 15

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
